{"js": "// Apply the edits described by the diff:\n// 1. \"...and added it to presentation.\" -> \"...and added it to sprint one presentation.\"\n// 2. \"Completed several other sections of the presentation.\" ->\n//    \"Completed several sections on the sprint one presentation.\"\n// 3. A new (empty) bold/Times-New-Roman paragraph is inserted right after the\n//    \"Completed several...\" bullet (before the existing blank paragraph that\n//    follows it).\n\nconst body = context.document.body;\n\n// ---- Change 1: hardware diagram bullet ----------------------------------\nconst bullet1 = body.search(\"Created hardware diagram that shows relationship between hardware necessary to our project and added it to presentation.\", { matchCase: true });\nbullet1.load(\"text\");\nawait context.sync();\n\nif (bullet1.items.length > 0) {\n  const target1 = bullet1.items[0].search(\"presentation.\", { matchCase: true });\n  target1.load(\"text\");\n  await context.sync();\n  target1.items[0].insertText(\"sprint one \", Word.InsertLocation.before);\n  await context.sync();\n}\n\n// ---- Change 2: \"Completed several other sections...\" bullet -------------\nconst bullet2 = body.search(\"Completed several other sections of the presentation.\", { matchCase: true });\nbullet2.load(\"text\");\nawait context.sync();\n\nif (bullet2.items.length > 0) {\n  const sentence = bullet2.items[0];\n\n  // Remove \"other \"\n  const otherHit = sentence.search(\"other \", { matchCase: true });\n  otherHit.load(\"text\");\n  await context.sync();\n  otherHit.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // \"of the\" -> \"on the\"\n  const ofTheHit = sentence.search(\"of the\", { matchCase: true });\n  ofTheHit.load(\"text\");\n  await context.sync();\n  ofTheHit.items[0].insertText(\"on the\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Insert \"sprint one \" before the final \"presentation.\"\n  const presHit = sentence.search(\"presentation.\", { matchCase: true });\n  presHit.load(\"text\");\n  await context.sync();\n  presHit.items[0].insertText(\"sprint one \", Word.InsertLocation.before);\n  await context.sync();\n}\n\n// ---- Change 3: insert a new blank paragraph after that bullet -----------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet bulletIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Completed several sections on the sprint one presentation.\") !== -1) {\n    bulletIndex = i;\n    break;\n  }\n}\n\nif (bulletIndex !== -1 && bulletIndex + 1 < paragraphs.items.length) {\n  const followingParagraph = paragraphs.items[bulletIndex + 1];\n  const newParagraph = followingParagraph.insertParagraph(\"\", Word.InsertLocation.before);\n  newParagraph.font.bold = true;\n  newParagraph.font.name = \"Times New Roman\";\n  newParagraph.font.size = 12;\n  newParagraph.spaceAfter = 0;\n  newParagraph.lineSpacing = 24;\n  await context.sync();\n}\n", "ps1": "# Apply the edits described by the diff:\n# 1. \"...and added it to presentation.\" -> \"...and added it to sprint one presentation.\"\n# 2. \"Completed several other sections of the presentation.\" ->\n#    \"Completed several sections on the sprint one presentation.\"\n# 3. A new (empty) bold/Times-New-Roman paragraph is inserted right after the\n#    \"Completed several...\" bullet (before the existing blank paragraph that\n#    follows it).\n\n$d = $word.ActiveDocument\n\n# ---- Change 1: hardware diagram bullet -----------------------------------\n$range1 = $d.Content\n$range1.Find.ClearFormatting()\n$range1.Find.Replacement.ClearFormatting()\n$found1 = $range1.Find.Execute(\n    \"Created hardware diagram that shows relationship between hardware necessary to our project and added it to presentation.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Created hardware diagram that shows relationship between hardware necessary to our project and added it to sprint one presentation.\",\n    2)\n\n# ---- Change 2: \"Completed several other sections...\" bullet -------------\n$range2 = $d.Content\n$range2.Find.ClearFormatting()\n$range2.Find.Replacement.ClearFormatting()\n$found2 = $range2.Find.Execute(\n    \"Completed several other sections of the presentation.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Completed several sections on the sprint one presentation.\",\n    2)\n\n# ---- Change 3: insert a new blank paragraph after that bullet -----------\n$paraCount = $d.Paragraphs.Count\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Contains(\"Completed several sections on the sprint one presentation.\")) {\n        $nextPara = $d.Paragraphs.Item($i + 1)\n        $nextPara.Range.InsertParagraphBefore()\n        break\n    }\n}\n"}
